$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = '71.809.32'
$ws.Cells.Item(2, 5).Value = '  +3.13%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '4.044.35'
$ws.Cells.Item(3, 5).Value = '  +2.82%  '

# Row 4
$ws.Cells.Item(4, 4).Value = '''0.999'
$ws.Cells.Item(4, 5).Value = '  +0.02%  '

# Row 5
$ws.Cells.Item(5, 4).Value = '''523.51'
$ws.Cells.Item(5, 5).Value = '  -2.07%  '

# Row 6
$ws.Cells.Item(6, 4).Value = '''148.88'
$ws.Cells.Item(6, 5).Value = '  +2.58%  '

# Row 7
$ws.Cells.Item(7, 5).Value = '  +1.05%  '

# Row 8
$ws.Cells.Item(8, 4).Value = '''0.999'
$ws.Cells.Item(8, 5).Value = '  +0.15%  '

# Row 9
$ws.Cells.Item(9, 4).Value = '''0.741'
$ws.Cells.Item(9, 5).Value = '  +1.65%  '

# Row 10
$ws.Cells.Item(10, 5).Value = '  +1.94%  '

# Row 11
$ws.Cells.Item(11, 4).Value = '''0.0000341'
$ws.Cells.Item(11, 5).Value = '  +0.24%  '

# Row 12
$ws.Cells.Item(12, 4).Value = '''46.70'
$ws.Cells.Item(12, 5).Value = '  +9.45%  '

# Row 13
$ws.Cells.Item(13, 4).Value = '''10.79'
$ws.Cells.Item(13, 5).Value = '  +3.68%  '

# Row 14
$ws.Cells.Item(14, 4).Value = '4.673.66'
$ws.Cells.Item(14, 5).Value = '  +2.49%  '

# Row 15
$ws.Cells.Item(15, 4).Value = '4.030.44'
$ws.Cells.Item(15, 5).Value = '  +2.20%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '''21.52'
$ws.Cells.Item(16, 5).Value = '  +8.43%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '''14.36'
$ws.Cells.Item(17, 5).Value = '  +2.35%  '

# Row 18
$ws.Cells.Item(18, 5).Value = '  +0.40%  '

# Row 19
$ws.Cells.Item(19, 5).Value = '  -1.72%  '

# Row 20
$ws.Cells.Item(20, 4).Value = '71.751.87'
$ws.Cells.Item(20, 5).Value = '  +3.21%  '

# Row 21
$ws.Cells.Item(21, 4).Value = '''442.04'
$ws.Cells.Item(21, 5).Value = '  +2.18%  '

# Row 22
$ws.Cells.Item(22, 5).Value = '  +5.61%  '

# Row 23
$ws.Cells.Item(23, 4).Value = '''94.60'
$ws.Cells.Item(23, 5).Value = '  +6.57%  '

# Row 24
$ws.Cells.Item(24, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(24, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(24, 4).Value = '''14.39'
$ws.Cells.Item(24, 5).Value = '  -0.81%  '

# Row 25
$ws.Cells.Item(25, 2).Value = 'RenderToken'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(25, 4).Value = '''12.36'
$ws.Cells.Item(25, 5).Value = '  +4.44%  '

# Row 26
$ws.Cells.Item(26, 4).Value = '''4.04'
$ws.Cells.Item(26, 5).Value = '  -2.46%  '

# Row 27
$ws.Cells.Item(27, 4).Value = '''11.17'
$ws.Cells.Item(27, 5).Value = '  +3.24%  '

# Row 28
$ws.Cells.Item(28, 4).Value = '''37.17'
$ws.Cells.Item(28, 5).Value = '  +1.34%  '

# Row 29
$ws.Cells.Item(29, 5).Value = '  +2.26%  '

# Row 30
$ws.Cells.Item(30, 4).Value = '''700.58'
$ws.Cells.Item(30, 5).Value = '  -0.26%  '

# Row 31
$ws.Cells.Item(31, 5).Value = '  +3.19%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +2.27%  '

# Row 33
$ws.Cells.Item(33, 4).Value = '''6.93'
$ws.Cells.Item(33, 5).Value = '  +13.28%  '

# Row 34
$ws.Cells.Item(34, 4).Value = '''67.69'
$ws.Cells.Item(34, 5).Value = '  -6.56%  '

# Row 35
$ws.Cells.Item(35, 4).Value = '0.0₃0912'
$ws.Cells.Item(35, 5).Value = '  +5.27%  '

# Row 36
$ws.Cells.Item(36, 4).Value = '''0.444'
$ws.Cells.Item(36, 5).Value = '  -2.64%  '

# Row 37
$ws.Cells.Item(37, 4).Value = '''40.83'
$ws.Cells.Item(37, 5).Value = '  +0.74%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  +5.75%  '

# Row 39
$ws.Cells.Item(39, 4).Value = '''3.54'
$ws.Cells.Item(39, 5).Value = '  +18.30%  '

# Row 40
$ws.Cells.Item(40, 5).Value = '  +0.33%  '

# Row 41
$ws.Cells.Item(41, 4).Value = '''0.998'
$ws.Cells.Item(41, 5).Value = '  -0.11%  '

# Row 42
$ws.Cells.Item(42, 4).Value = '''0.0490'
$ws.Cells.Item(42, 5).Value = '  +1.52%  '

# Row 43
$ws.Cells.Item(43, 4).Value = '''3.13'
$ws.Cells.Item(43, 5).Value = '  +0.63%  '

# Row 44
$ws.Cells.Item(44, 5).Value = '  -0.09%  '

# Row 45
$ws.Cells.Item(45, 4).Value = '''3.53'
$ws.Cells.Item(45, 5).Value = '  +4.04%  '

# Row 46
$ws.Cells.Item(46, 2).Value = 'Stellar'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(46, 4).Value = '''0.146'
$ws.Cells.Item(46, 5).Value = '  +2.57%  '

# Row 47
$ws.Cells.Item(47, 2).Value = 'Stacks'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(47, 4).Value = '''3.21'
$ws.Cells.Item(47, 5).Value = '  +0.18%  '

# Row 48
$ws.Cells.Item(48, 4).Value = '''0.000283'
$ws.Cells.Item(48, 5).Value = '  +19.77%  '

# Row 49
$ws.Cells.Item(49, 4).Value = '''9.19'
$ws.Cells.Item(49, 5).Value = '  +5.70%  '

# Row 50
$ws.Cells.Item(50, 4).Value = '''3.38'
$ws.Cells.Item(50, 5).Value = '  +1.44%  '

# Row 51
$ws.Cells.Item(51, 4).Value = '0.0₆0347'
$ws.Cells.Item(51, 5).Value = '  -2.40%  '
